$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2236.182
$ws.Range("I19").Value = 2426.8
$ws.Range("J19").Value = 2077.3333
$ws.Range("K19").Value = 2426.8
$ws.Range("L19").Value = 2077.3333
$ws.Range("M19").Value = -2251.8
$ws.Range("N19").Value = -2427.3333
$ws.Range("H53").Value = 174.66667
$ws.Range("I53").Value = 189.9375
$ws.Range("J53").Value = 157.21428
$ws.Range("K53").Value = 189.9375
$ws.Range("L53").Value = 157.21428
$ws.Range("M53").Value = 447.0625
$ws.Range("N53").Value = -1431.21428
$ws.Range("H80").Value = 1088.5
$ws.Range("I80").Value = 1241.8
$ws.Range("J80").Value = 833
$ws.Range("K80").Value = 3725.4
$ws.Range("L80").Value = 2499
$ws.Range("M80").Value = -2727.4
$ws.Range("N80").Value = -4495
$ws.Range("H83").Value = 1088.5
$ws.Range("I83").Value = 1241.8
$ws.Range("J83").Value = 833
$ws.Range("K83").Value = 11176.2
$ws.Range("L83").Value = 7497
$ws.Range("M83").Value = -6184.199999999999
$ws.Range("N83").Value = -17481

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1787
$ws.Range("I2").Value = 1853.7273
$ws.Range("J2").Value = 1664.6666
$ws.Range("K2").Value = 1853.7273
$ws.Range("L2").Value = 1664.6666
$ws.Range("M2").Value = -1740.7273
$ws.Range("N2").Value = -1890.6666
$ws.Range("H45").Value = 11726.272
$ws.Range("I45").Value = 11459.637
$ws.Range("J45").Value = 12259.546
$ws.Range("K45").Value = 11459.637
$ws.Range("L45").Value = 12259.546
$ws.Range("M45").Value = -11082.637
$ws.Range("N45").Value = -13013.546
$ws.Range("H46").Value = 6871.7144
$ws.Range("I46").Value = 5000
$ws.Range("J46").Value = 7183.6665
$ws.Range("K46").Value = 5000
$ws.Range("L46").Value = 7183.6665
$ws.Range("M46").Value = -4681
$ws.Range("H61").Value = 1735.7916
$ws.Range("I61").Value = 1669.6364
$ws.Range("J61").Value = 2463.5
$ws.Range("K61").Value = 1669.6364
$ws.Range("L61").Value = 2463.5
$ws.Range("M61").Value = -1457.6364
$ws.Range("H74").Value = 1545
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 1545
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 1545
$ws.Range("N74").Value = -3293
$ws.Range("H77").Value = 1545
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 1545
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 7725
$ws.Range("N77").Value = -16461
$ws.Range("H80").Value = 18400
$ws.Range("I80").Value = 12000
$ws.Range("J80").Value = 20000
$ws.Range("K80").Value = 12000
$ws.Range("L80").Value = 20000
$ws.Range("M80").Value = -11002
$ws.Range("N80").Value = -21996
$ws.Range("H83").Value = 18400
$ws.Range("I83").Value = 12000
$ws.Range("J83").Value = 20000
$ws.Range("K83").Value = 36000
$ws.Range("L83").Value = 60000
$ws.Range("M83").Value = -31008
$ws.Range("N83").Value = -69984
$ws.Range("H97").Value = 1063.7241
$ws.Range("I97").Value = 716.6667
$ws.Range("J97").Value = 5749
$ws.Range("K97").Value = 716.6667
$ws.Range("L97").Value = 5749
$ws.Range("M97").Value = -220.6667
$ws.Range("H102").Value = 2435.4285
$ws.Range("I102").Value = 2709.8
$ws.Range("J102").Value = 1749.5
$ws.Range("K102").Value = 2709.8
$ws.Range("L102").Value = 1749.5
$ws.Range("M102").Value = -1087.8
$ws.Range("H116").Value = 1787
$ws.Range("I116").Value = 1853.7273
$ws.Range("J116").Value = 1664.6666
$ws.Range("K116").Value = 1853.7273
$ws.Range("L116").Value = 1664.6666
$ws.Range("M116").Value = 440.2727
$ws.Range("N116").Value = -6252.6666
$ws.Range("H132").Value = 3026.6553
$ws.Range("I132").Value = 2893.0588
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8679.1764
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6149.1764
$ws.Range("H136").Value = 1735.7916
$ws.Range("I136").Value = 1669.6364
$ws.Range("J136").Value = 2463.5
$ws.Range("K136").Value = 5008.9092
$ws.Range("L136").Value = 7390.5
$ws.Range("M136").Value = -2458.9092
$ws.Range("M74").ClearContents()
$ws.Range("M77").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1787
$ws.Range("I3").Value = 1853.7273
$ws.Range("J3").Value = 1664.6666
$ws.Range("K3").Value = 1853.7273
$ws.Range("L3").Value = 1664.6666
$ws.Range("M3").Value = -1739.7273
$ws.Range("N3").Value = -1892.6666
$ws.Range("H105").Value = 12128.792
$ws.Range("I105").Value = 9013.857
$ws.Range("J105").Value = 33933.332
$ws.Range("K105").Value = 9013.857
$ws.Range("L105").Value = 33933.332
$ws.Range("M105").Value = -7266.857
$ws.Range("H107").Value = 3747.1428
$ws.Range("I107").Value = 2651.5715
$ws.Range("J107").Value = 4842.7144
$ws.Range("K107").Value = 2651.5715
$ws.Range("L107").Value = 4842.7144
$ws.Range("M107").Value = -731.5715
$ws.Range("N107").Value = -8682.714400000001
$ws.Range("H134").Value = 2823.7646
$ws.Range("I134").Value = 2049.6099
$ws.Range("J134").Value = 5997.8
$ws.Range("K134").Value = 6148.8297
$ws.Range("L134").Value = 17993.4
$ws.Range("M134").Value = -3613.8297
$ws.Range("N134").Value = -23063.4

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4818.85
$ws.Range("I31").Value = 2131.5715
$ws.Range("J31").Value = 6265.846
$ws.Range("K31").Value = 2131.5715
$ws.Range("L31").Value = 6265.846
$ws.Range("M31").Value = -1836.5715
$ws.Range("N31").Value = -6855.846
$ws.Range("H34").Value = 4818.85
$ws.Range("I34").Value = 2131.5715
$ws.Range("J34").Value = 6265.846
$ws.Range("K34").Value = 2131.5715
$ws.Range("L34").Value = 6265.846
$ws.Range("M34").Value = -1929.5715
$ws.Range("N34").Value = -6669.846
$ws.Range("H105").Value = 2392.5715
$ws.Range("I105").Value = 2710
$ws.Range("J105").Value = 1599
$ws.Range("K105").Value = 2710
$ws.Range("L105").Value = 1599
$ws.Range("M105").Value = -963
$ws.Range("N105").Value = -5093
$ws.Range("H107").Value = 1169.2894
$ws.Range("I107").Value = 1078.76
$ws.Range("J107").Value = 1343.3846
$ws.Range("K107").Value = 1078.76
$ws.Range("L107").Value = 1343.3846
$ws.Range("M107").Value = 841.24

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 240.66667
$ws.Range("I23").Value = 199.8
$ws.Range("J23").Value = 261.1
$ws.Range("K23").Value = 599.4000000000001
$ws.Range("L23").Value = 783.3000000000001
$ws.Range("M23").Value = -364.4000000000001
$ws.Range("N23").Value = -1253.3
$ws.Range("H39").Value = 47452.273
$ws.Range("I39").Value = 67500
$ws.Range("J39").Value = 4492.857
$ws.Range("K39").Value = 202500
$ws.Range("L39").Value = 13478.571
$ws.Range("M39").Value = -202206
$ws.Range("N39").Value = -14066.571
$ws.Range("H55").Value = 7416264.5
$ws.Range("I55").Value = 145357
$ws.Range("J55").Value = 12505900
$ws.Range("K55").Value = 436071
$ws.Range("L55").Value = 37517700
$ws.Range("M55").Value = -435894
$ws.Range("N55").Value = -37518054
$ws.Range("H59").Value = 2252.5
$ws.Range("I59").Value = 1605
$ws.Range("J59").Value = 2900
$ws.Range("K59").Value = 4815
$ws.Range("L59").Value = 8700
$ws.Range("M59").Value = -4275
$ws.Range("H68").Value = 2429.7222
$ws.Range("I68").Value = 2680.8333
$ws.Range("J68").Value = 2304.1667
$ws.Range("K68").Value = 8042.499899999999
$ws.Range("L68").Value = 6912.500100000001
$ws.Range("M68").Value = -7231.499899999999
$ws.Range("N68").Value = -8534.500100000001
$ws.Range("H71").Value = 2429.7222
$ws.Range("I71").Value = 2680.8333
$ws.Range("J71").Value = 2304.1667
$ws.Range("K71").Value = 24127.4997
$ws.Range("L71").Value = 20737.5003
$ws.Range("M71").Value = -20071.4997
$ws.Range("N71").Value = -28849.5003

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2296.7693
$ws.Range("I102").Value = 2296.125
$ws.Range("J102").Value = 2304.5
$ws.Range("K102").Value = 2296.125
$ws.Range("L102").Value = 2304.5
$ws.Range("M102").Value = -674.125
$ws.Range("N102").Value = -5548.5
$ws.Range("H132").Value = 4165.61
$ws.Range("I132").Value = 4059.5278
$ws.Range("J132").Value = 4929.4
$ws.Range("K132").Value = 12178.5834
$ws.Range("L132").Value = 14788.2
$ws.Range("M132").Value = -9648.5834
$ws.Range("N132").Value = -19848.2

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1068.4
$ws.Range("I46").Value = 1068.4
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1068.4
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -880.4000000000001
$ws.Range("H82").Value = 812.2353000000001
$ws.Range("I82").Value = 647.53845
$ws.Range("J82").Value = 1347.5
$ws.Range("K82").Value = 647.53845
$ws.Range("L82").Value = 1347.5
$ws.Range("M82").Value = -286.53845
$ws.Range("N82").Value = -2069.5
$ws.Range("H85").Value = 812.2353000000001
$ws.Range("I85").Value = 647.53845
$ws.Range("J85").Value = 1347.5
$ws.Range("K85").Value = 647.53845
$ws.Range("L85").Value = 1347.5
$ws.Range("M85").Value = 600.46155
$ws.Range("N85").Value = -3843.5
$ws.Range("N46").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("H97").Value = 18078.5
$ws.Range("I97").Value = 29900
$ws.Range("J97").Value = 15714.2
$ws.Range("K97").Value = 29900
$ws.Range("L97").Value = 15714.2
$ws.Range("M97").Value = -28909
$ws.Range("N97").Value = -17696.2
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()
